$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (small)
$ws.Range("B2").Value = 6.692
$ws.Range("C2").Value = 47.671
$ws.Range("D2").Value = 1.088
$ws.Range("E2").Value = 1.939
$ws.Range("F2").Value = 57.39

# Row 3 (medium)
$ws.Range("B3").Value = 10.925
$ws.Range("C3").Value = 12.272
$ws.Range("D3").Value = 1.04
$ws.Range("E3").Value = 0.07099999999999999
$ws.Range("F3").Value = 24.308

# Row 4 (large)
$ws.Range("B4").Value = 5.699
$ws.Range("C4").Value = 6.834
$ws.Range("D4").Value = 0.355
$ws.Range("E4").Value = 0.024
$ws.Range("F4").Value = 12.912

# Row 5 (huge)
$ws.Range("B5").Value = 0.236
$ws.Range("C5").Value = 0.047
$ws.Range("F5").Value = 0.283

# Row 6 (unknown_sz)
$ws.Range("B6").Value = 0.497
$ws.Range("C6").Value = 3.925
$ws.Range("D6").Value = 0.118
$ws.Range("E6").Value = 0.5679999999999999
$ws.Range("F6").Value = 5.108

# Row 7 (COL_TOT)
$ws.Range("B7").Value = 24.049
$ws.Range("C7").Value = 70.749
$ws.Range("D7").Value = 2.601
$ws.Range("E7").Value = 2.602
$ws.Range("F7").Value = 100.001
